$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# trafo_id -> gridnode_id refactor: rename the header in column J
$ws.Range("J1").Value = "gridnode_id"

# update the active selection to match the authored workbook state
[void]$ws.Range("G8").Select()
